$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date serial numbers; rows 2-12 all had the
# same "last changed" date (serial 45178 = 2023-09-09) which is bumped by
# one day (serial 45179 = 2023-09-10).
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45179
}
